{"js": "// Replace \"overlay mesh\" with \"Docker ingress\" in the proxy port-ranges\n// paragraph, and move the \"_GoBack\" bookmark from its old location (end of\n// the \"node_poll_seconds ... 30 seconds.\" paragraph) to just after the new\n// \"Docker ingress\" text, matching how Word leaves its last-edit-position\n// bookmark behind after making this change.\n\n// 1) Drop the old \"_GoBack\" bookmark wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the sentence that needs the wording change.\nconst sentence =\n  \"NeonClusters reserves a block of 100 ports on the overlay mesh network \" +\n  \"for each of the public and private proxies.\";\nconst searchResults = context.document.body.search(sentence, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the \"overlay mesh\" sentence to update.');\n}\n\nconst target = searchResults.items[0];\n\n// 3) Replace the sentence with OOXML that already has the new wording split\n//    into separate runs around the re-inserted \"_GoBack\" bookmark (this is\n//    what Word itself produces when it records the edit position).\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">NeonClusters reserves a block of 100 ports on the </w:t></w:r><w:r><w:t>Docker ingress</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\"> network for each of the public and private proxies.</w:t></w:r></w:p></w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace \"overlay mesh\" with \"Docker ingress\" in the proxy port-ranges\n# paragraph, and move the \"_GoBack\" bookmark from its old location (end of\n# the \"node_poll_seconds ... 30 seconds.\" paragraph) to just after the new\n# \"Docker ingress\" text, matching how Word leaves its last-edit-position\n# bookmark behind after making this change.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Find the sentence that needs the wording change.\n$sentence = \"NeonClusters reserves a block of 100 ports on the overlay mesh network for each of the public and private proxies.\"\n$rng = $d.Content\n$rng.Find.Execute($sentence)\nif (-not $rng.Find.Found) {\n    throw 'Could not find the \"overlay mesh\" sentence to update.'\n}\n\n# Clear the matched range so the OOXML below fully replaces it.\n$rng.Text = \"\"\n\n# 3) Insert OOXML that has the new wording split into separate runs around\n#    the re-inserted \"_GoBack\" bookmark (this is what Word itself produces\n#    when it records the edit position).\n$newXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p w:rsidR=\"00F537C8\" w:rsidRDefault=\"00F537C8\" w:rsidP=\"00F537C8\"><w:r><w:t xml:space=\"preserve\">NeonClusters reserves a block of 100 ports on the </w:t></w:r><w:r><w:t>Docker ingress</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\"> network for each of the public and private proxies.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$rng.InsertXML($newXml)\n"}
